$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "22"
$ws.Range("M2").Style = "Normal"

$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "19"
$ws.Range("M4").Style = "Normal"

$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = "11"
$ws.Range("M5").Style = "Normal"

$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "6"
$ws.Range("M6").Style = "Normal"

$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = "6"
$ws.Range("M7").Style = "Normal"

$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = "5"
$ws.Range("M8").Style = "Normal"

$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "3"
$ws.Range("M9").Style = "Normal"
